# Clear stale accomplishment-tracking columns now that the latest
# status/accomplishment files (as of May) have been incorporated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# L2 (PROJECT ALLOCATION for the first data row) is no longer populated.
$ws.Range("L2").ClearContents()

# Columns AB:AK (site/building reverted/not-started/procurement/ongoing/
# completed breakdown) and AM (DIFFERENCE) are dropped for every data row
# (2-28); column AL (PREVIOUS ACCOMPLISHMENT) is left untouched.
$ws.Range("AB2:AK28").ClearContents()
$ws.Range("AM2:AM28").ClearContents()
